# Update "想去人数" (F column) counts on the "展览" sheet and the
# combined "全部类型" sheet to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 33
$ws1.Range("F3").Value = 303
$ws1.Range("F5").Value = 2607
$ws1.Range("F6").Value = 1868
$ws1.Range("F7").Value = 357
$ws1.Range("F8").Value = 111
$ws1.Range("F9").Value = 917
$ws1.Range("F10").Value = 180

# --- Sheet "全部类型" (all types, combines 展览 + 演出) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 33
$ws4.Range("F3").Value = 303
$ws4.Range("F5").Value = 2607
$ws4.Range("F6").Value = 1868
$ws4.Range("F7").Value = 357
$ws4.Range("F9").Value = 111
$ws4.Range("F10").Value = 917
$ws4.Range("F11").Value = 180
